$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(22, 38.35754998722376, -16.75077860103589, 9.333282041253318),
    @(23, 4.113532688307743, -10.04574092130529, 6.603543145316186),
    @(24, -28.9259203975992, 2.342653398928785, -14.07138639651451),
    @(25, 42.44476240762218, -44.82572951820286, -4.690105651476593),
    @(26, 47.95408476361651, -34.91908646222309, 8.182450048671846),
    @(27, 16.99495727823363, -33.59686962269733, 18.0336912966659),
    @(28, 0.5081373297654075, -34.04575231800982, 16.78914682761511),
    @(29, -27.83808203513578, 3.896602891246781, -12.58596580072954),
    @(30, 21.10497084907895, -35.42225003686769, -4.233091697929673),
    @(31, 45.68115871737581, -33.5533879203098, 10.18618937901121)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
